$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.931.43"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "2.376.73"
$ws.Range("E3").Value = "  +2.93%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "302.44"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "98.78"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.94%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.571"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.29%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "34.39"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.82%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0803"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "7.19"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "2.737.60"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").Value = "2.374.79"
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "13.77"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.813"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "45.882.60"
$ws.Range("E18").Value = "  -1.18%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.96"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "0.0₃0978"
$ws.Range("E20").Value = "  +4.25%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.08"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.04%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "67.03"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "245.16"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.26%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.85"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("E25").Value = "  +0.23%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.94"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -2.67%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "39.21"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -6.96%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.21"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -3.22%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.80"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "21.40"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +6.63%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.79"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +21.84%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.79"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +7.26%  "
$ws.Range("E33").Value = "  -2.68%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "145.69"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.90%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.0777"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.92%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.114"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +2.06%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.90"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +5.53%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.96"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "14.93"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.49%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.0302"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E42").Value = "  -6.18%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "1.871.31"
$ws.Range("E44").Value = "  +1.74%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "90.79"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("E46").Value = "  -11.40%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "8.39"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +5.81%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "15.02"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +8.79%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.187"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -5.36%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.608.74"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "97.92"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.40%  "
